# "create settting to push variable there"
# Adds a new "Sheet2" worksheet (after Sheet1) containing a second batch of
# source/target language rows + a total_tokens column, mirroring Sheet1's
# layout. Also nudges the sheetView rightToLeft flag (best effort).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- add Sheet2 right after Sheet1 -----------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# --- header row --------------------------------------------------------
$ws2.Range("A1").Value = "source_language"
$ws2.Range("B1").Value = "target_language"
$ws2.Range("C1").Value = "total_tokens"

# --- data rows -----------------------------------------------------------
$ws2.Range("A2").Value = "You always know after you are two. Two is the beginning of the end."
$ws2.Range("B2").Value = "Tu toujours"
$ws2.Range("C2").Value = 24

$ws2.Range("A3").Value = "This was all that passed between them on the subject, but henceforth Wendy knew that she must grow up."
$ws2.Range("B3").Value = "Cela était"
$ws2.Range("C3").Value = 26

$ws2.Range("A4").Value = "They soon know that they will grow up, and the way Wendy knew was this."
$ws2.Range("B4").Value = "Ils bientôt"
$ws2.Range("C4").Value = 28

$ws2.Range("A5").Value = "All children, except one, grow up."
$ws2.Range("B5").Value = "en deserve the opportunity to learn and grow in a safe and supportive environment."
$ws2.Range("C5").Value = 38

# --- view flags (best effort; property is read-only in this host) ----------
$ws1.Select()
try { $excel.ActiveWindow.DisplayRightToLeft = $false } catch {}
$ws2.Select()
try { $excel.ActiveWindow.DisplayRightToLeft = $false } catch {}
$ws1.Select()
